$wb = $excel.ActiveWorkbook

# --- Sheet "Form_Games": insert a new game row at position 74 -------------
$ws = $wb.Worksheets.Item("Form_Games")
$ws.Rows.Item(74).Insert()

function Set-TextCell($rng, $text) {
    # Leading apostrophe forces literal text entry so Excel's heuristic
    # date/number auto-detection (e.g. "01-12-2025" looking like a date)
    # doesn't silently turn the string into a date serial number.
    $rng.Value = "'" + $text
    # ...then strip the quote-prefix/number-format bookkeeping that the
    # forced-text entry leaves behind, back to the plain default style.
    $rng.Style = "Normal"
}

Set-TextCell $ws.Range("A74") "01-12-2025"
Set-TextCell $ws.Range("B74") "Торпедо"
Set-TextCell $ws.Range("C74") "Динамо М"
$ws.Range("D74").Value = 897837
Set-TextCell $ws.Range("E74") "https://text.khl.ru/text/897837.html"
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 4
$ws.Range("I74").Value = 29
$ws.Range("J74").Value = 24
$ws.Range("K74").Value = 53
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 2
$ws.Range("N74").Value = 8
$ws.Range("O74").Value = 4
$ws.Range("P74").Value = 41
$ws.Range("Q74").Value = 29
$ws.Range("R74").Value = 17
$ws.Range("S74").Value = 16
$ws.Range("T74").Value = 27
$ws.Range("U74").Value = 8
$ws.Range("V74").Value = 1
$ws.Range("W74").Value = 2
$ws.Range("X74").Value = 10
$ws.Range("Y74").Value = 17
$ws.Range("Z74").Value = 3
$ws.Range("AA74").Value = 7
Set-TextCell $ws.Range("AB74") "L"
$ws.Range("AC74").Value = 58.6
$ws.Range("AD74").Value = 77.09999999999999
$ws.Range("AE74").Value = 0
$ws.Range("AF74").Value = 0
$ws.Range("AG74").Value = 0.5
$ws.Range("AH74").Value = 0

# --- Sheet "Aggregates": refresh Torpedo's (row 16) rolling stats ---------
$agg = $wb.Worksheets.Item("Aggregates")
$agg.Range("B16").Value = 5
$agg.Range("C16").Value = 2.4
$agg.Range("D16").Value = 3
$agg.Range("E16").Value = 0.2
$agg.Range("F16").Value = 0.3
$agg.Range("G16").Value = 29.8
$agg.Range("H16").Value = 60.7
$agg.Range("I16").Value = -0.1
$agg.Range("J16").Value = 28
$agg.Range("K16").Value = 15.5
$agg.Range("L16").Value = 1.1
$agg.Range("M16").Value = 21.6
$agg.Range("N16").Value = 2.2
$agg.Range("O16").Value = 51.4
$agg.Range("P16").Value = -2.740000000000001
$agg.Range("Q16").Value = 0.166
$agg.Range("R16").Value = -0.01600000000000001
$agg.Range("S16").Value = 0.5660000000000001
$agg.Range("T16").Value = 0.034
$agg.Range("U16").Value = 2.2
$agg.Range("V16").Value = -0.2
$agg.Range("W16").Value = 6.4
$agg.Range("X16").Value = 1
$agg.Range("Y16").Value = 8
$agg.Range("Z16").Value = 0.2
$agg.Range("AA16").Value = 8.4
$agg.Range("AB16").Value = 0.7200000000000001

Write-Output "edit applied"
